# Update Name of Algo
# Apply the updated RandomForest result values to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -6.516899999999993
$ws.Range("C3").Value = -11.12149999999999
$ws.Range("D5").Value = -8.648799999999992
$ws.Range("C14").Value = -12.5019
$ws.Range("C21").Value = -13.19350000000001
$ws.Range("C23").Value = -12.27120000000001
$ws.Range("C25").Value = -11.0584

$wb.Save()
